# Added Both ON/OFF vendors in Download SF List
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell (row 1) - bold/centered like the other header cells
$ws.Cells.Item(1, 34).Value = "On/Off"
$ws.Cells.Item(1, 34).Style = $ws.Cells.Item(1, 33).Style

# New placeholder cell (row 2) - plain style, matches AB2 (no explicit style)
$ws.Cells.Item(2, 34).Value = "{vendor:on_off_status}"

# Match column AH width to the rest of the template (bestFit column)
$ws.Columns.Item(34).ColumnWidth = 19.3

# Update the view: scroll to show the new column, move the active selection
$ws.Application.ActiveWindow.ScrollColumn = 28
$ws.Range("AI7").Select()
